$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EoBSDwEC")

# New rows for additional fuels (added first so shared-string indices land before the header text)
$ws.Range("A8").Value = "kerosene"
$ws.Range("B8").Value = -0.15
$ws.Range("C8").Value = -0.15
$ws.Range("D8").Value = -0.25

$ws.Range("A9").Value = "heavy or residual fuel oil"
$ws.Range("B9").Value = -0.15
$ws.Range("C9").Value = -0.15
$ws.Range("D9").Value = -0.25

$ws.Range("A10").Value = "LPG propane or butane"
$ws.Range("B10").Value = -0.15
$ws.Range("C10").Value = -0.15
$ws.Range("D10").Value = -0.25

$ws.Range("A11").Value = "hydrogen"
$ws.Range("B11").Value = -0.15
$ws.Range("C11").Value = -0.15
$ws.Range("D11").Value = -0.25

# Update header row
$ws.Range("A1").Value = "Elasticity by Fuel (dimensionless)"
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Column width tweaks (closest achievable approximations of the saved widths)
$ws.Columns.Item(1).ColumnWidth = 24.166666666666668
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 19.0
$ws.Columns.Item(4).ColumnWidth = 13.333333333333332
